# Update average_county_temperature (column AD) for rows 9-36 with
# the updated NOAA temperature data value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValue = 19.30324074074072

for ($row = 9; $row -le 36; $row++) {
    $ws.Range("AD$row").Value = $newValue
}
